$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 01:05"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1384033
$ws.Range("C4").Value = 16395
$ws.Range("D4").Value = 260355
$ws.Range("E4").Value = 1041975
$ws.Range("F4").Value = 16481
$ws.Range("G4").Value = 916
$ws.Range("H4").Value = 81703

# Row 51 - Chequia
$ws.Range("B51").Value = 8176
$ws.Range("C51").Value = 53
$ws.Range("D51").Value = 4711
$ws.Range("E51").Value = 3183
$ws.Range("F51").Value = 40
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 282

# Row 92 - Republica de Yibuti
$ws.Range("B92").Value = 1227
$ws.Range("C92").Value = 17
$ws.Range("D92").Value = 872
$ws.Range("E92").Value = 352

# Row 107 - Niger
$ws.Range("B107").Value = 832
$ws.Range("C107").Value = 11
$ws.Range("D107").Value = 637
$ws.Range("E107").Value = 149

# Row 115 - Uruguay
$ws.Range("B115").Value = 711
$ws.Range("C115").Value = 4
$ws.Range("D115").Value = 523
$ws.Range("E115").Value = 169

$wb.Save()
